# Activity_Utilization.xlsx — "completed happy path scenarios for the utilization feature"
#
# The three API-scenario sheets were mis-named (sessions/summary/last_login
# were scrambled relative to their actual content). Rename them correctly,
# then fill in the previously-blank "Get random equipmentId" pre-requisite
# column and the HTTP method/URI/description cells that make the
# equipment-session and equipment-summary happy-path rows complete.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # GET_equipment                (untouched)
$ws2 = $wb.Worksheets.Item(2)   # was GET_equipId_sessions
$ws3 = $wb.Worksheets.Item(3)   # was GET_equipId_summary
$ws4 = $wb.Worksheets.Item(4)   # was GET_equipId_last_login

# --- Fix the sheet names so they match what each tab actually tests ------
$ws2.Name = "GET_last_login"
$ws3.Name = "GET_equipment_session"
$ws4.Name = "GET_equipment_summary"

# --- Complete the happy-path rows -----------------------------------------
# GET_equipment_session (ws3): the session endpoint's URI, then ...
$ws3.Range("F2").Value = "/activity/v1/equipment/{equipmentId}/sessions?pageSize=50&startTimestamp=2021-5-1&endTimestamp=2022-5-1"
# GET_equipment_summary (ws4): the summary endpoint's URI, then ...
$ws4.Range("F2").Value = "/activity/v1/equipment/{equipmentId}/summary?pageSize=50&startTimestamp=2021-5-1&endTimestamp=2022-5-1"
# ... the session sheet's description, and both sheets' "Get random
# equipmentId" pre-requisite (already used elsewhere, so no new string).
$ws3.Range("C2").Value = "Get equipment session"
$ws3.Range("D2").Value = "Get random equipmentId"
$ws4.Range("D2").Value = "Get random equipmentId"

# Narrow column D on the session sheet now that it holds shorter text.
$ws3.Columns.Item(4).ColumnWidth = 28.8

# --- Update selections / active tab ---------------------------------------
# Leave ws1 exactly as-is. Move each sheet's selection to D2, then finish
# with ws4 active so it becomes the selected tab (matches activeTab="3").
$ws2.Range("D2").Select()
$ws3.Range("D2").Select()
$ws4.Range("D2").Select()
